# Adds the "Hoja1" worksheet (expanded per-base variable dictionary) after "Variables",
# matching the "Anadir variable de delitos cortos" commit.

$wb = $excel.ActiveWorkbook

# --- 1. Set the selection on "Variables" before switching tabs (matches the
#        diff's sheet1 sheetView: selection sqref="A1:D1", no more tabSelected/topLeftCell) ---
$wsVariables = $wb.Worksheets.Item(1)
$wsVariables.Range("A1:D1").Select()

# --- 2. Insert the new worksheet right after "Variables" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsHoja1 = $wb.Worksheets.Add($null, $lastSheet)
$wsHoja1.Name = "Hoja1"

# --- 3. Header row (row 1): Nombre / Descripcion / Tipo / Base de origen ---
$wsHoja1.Cells.Item(1,1).Value = "Nombre"
$wsHoja1.Cells.Item(1,2).Value = "Descripción"
$wsHoja1.Cells.Item(1,3).Value = "Tipo "
$wsHoja1.Cells.Item(1,4).Value = "Base de origen "

# B1:D1 reuse the bold/centered header look already used on sheet "Variables"
$hdrRest = $wsHoja1.Range("B1:D1")
$hdrRest.Font.Bold = $true
$hdrRest.Font.Size = 11
$hdrRest.HorizontalAlignment = -4108

# A1 gets its own bold header look (smaller, 9pt)
$hdrA1 = $wsHoja1.Range("A1")
$hdrA1.Font.Bold = $true
$hdrA1.Font.Size = 9
$hdrA1.HorizontalAlignment = -4108

# --- 4. Data rows: column A, rows 2-75 ---
$rows = @'
2	id_exp	2
3	id_per_acusada	2
4	aparece_en_bases	2
5	base_asuntos	2
6	base_sitjurid	2
7	base_sol_alternas	2
8	base_medida_cautelar	2
9	base_sentencias	2
10	num_alcaldias_asunto	2
11	num_consignacion_asunto	2
12	num_comision_asunto	2
13	num_realizacion_asunto	2
14	num_delitos_asuntos	2
15	num_alcaldias_sitjurid	2
16	num_consignacion_sitjurid	2
17	num_comision_sitjurid	2
18	num_realizacion_sitjurid	2
19	num_delitos_sitjurid	2
20	num_alcaldias_alternas	4
21	num_consignacion_alternas	2
22	num_comision_alternas	2
23	num_realizacion_alternas	2
24	num_delitos_alternas	2
25	num_alcaldias_cautelares	2
26	num_consignacion_cautelares	2
27	num_comision_cautelares	2
28	num_realizacion_cautelares	2
29	num_delitos_cautelares	2
30	num_medidas_cautelares	2
31	num_alcaldias_sentencia	2
32	num_consignacion_sentencia	2
33	num_comision_sentencia	2
34	num_realizacion_sentencia	2
35	num_delitos_sentencia	2
36	num_terminacion	2
37	num_ppo	2
38	year_asunto	2
39	year_sitjurid	2
40	year_alternas	2
41	year_cautelares	2
42	year_sentencia	2
43	month_asunto	2
44	month_sitjurid	2
45	month_alternas	2
46	month_cautelares	2
47	month_sentencia	2
48	date_asunto	2
49	date_sijurid	2
50	date_alternas	4
51	date_cautelares	2
52	date_sentencia	2
53	sexo_acusada_asunto	2
54	sexo_acusada_sitjurid	2
55	sexo_acusada_alternas	2
56	sexo_acusada_cautelares	2
57	sexo_sentenciada	2
58	edad_acusada_asunto	2
59	edad_acusada_sitjurid	2
60	edad_acusada_alternas	2
61	edad_acusada_cautelares	2
62	edad_sentenciada	2
63	con_terminacion	2
64	tipo_terminacion	2
65	con_ppo	2
66	tipo_ppo	2
67	materia_asunto	2
68	materia_sitjurid	2
69	materia_alternas	2
70	materia_cautelares	2
71	materia_sentencia	2
72	c_con_detenido_asunto	2
73	c_sin_detenido_asunto	2
74	c_culposo_asunto	2
75	c_doloso_asunto	2
'@ -split "`n"

foreach ($line in $rows) {
    if ([string]::IsNullOrWhiteSpace($line)) { continue }
    $parts = $line -split "`t"
    $r = [int]$parts[0]
    $text = $parts[1]
    $styleCode = $parts[2]
    $cell = $wsHoja1.Cells.Item($r, 1)
    $cell.Value = $text
    if ($styleCode -eq "4") {
        $cell.Font.Size = 9
    } else {
        $cell.Font.Name = "Lucida Console"
        $cell.Font.Size = 9
        $cell.Font.Family = 3
        $cell.VerticalAlignment = -4108
    }
}

# --- 5. Column widths ---
$wsHoja1.Columns.Item(1).ColumnWidth = 5.29
$wsHoja1.Columns.Item(4).ColumnWidth = 14.43

# --- 6. View state: Hoja1 becomes the active/visible tab, scrolled near A31,
#        with A30 selected (matches the diff's sheet2 sheetView) ---
$wsHoja1.Activate()
$wsHoja1.Range("A30").Select()

Write-Host "Hoja1 added with $($rows.Count) data rows"
